$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "best camera and display and best performance in the  segment ..... excellent for gaming... pubg experiece is so perfect and the g95 processor gives a best performance"
$ws.Range("B3").Value = "When we talk about a brand loyal person....so here i am....it is pleasure to say that...you are coming to a very superb brand which was already known by Realme....this Realme 7 in mist blue colour is awesome as well as their performance is un-beatable....i mean wow yr.....this hand-set is having a different personality when it was carried whether inside the house or outside in a party or wedding function....what a wonderful performance bro....💓So i prerer you to purchase this phone....and i..."
$ws.Range("B4").Value = "Writting this after using for 1 month.Pros.1. Camera is too good at this price range. 2. Compact design (Easily usable in one hand).3. No heavy bloatware installed other than some realme apps.4. Battery durability is quite good.5. No heating issue found even after playing PUBG for 3 hours.Cons.1. Sometimes photos are over exposed. Camera software optimization needs to be improved.2. Sound quality is disappointing3. Display brightness is not enough for a sunny day. You'll get trouble..."
$ws.Range("B5").Value = "Best phone in the range Camera is better ✔️Battery backup is ok 👍Fast charging 💯Look is mustPerfomence good ⛹️Aal is very good"
$ws.Range("B6").Value = "Realme7 It's a branded mobile To value of money and good battery backup Display was awesomeSound quality was nice Good gaming 👍 Camera was extremely super good quality pics Side fingerprint sensor was NYC and speed front camera super Worth for maney No lag no struck 30watts charger was super"
$ws.Range("B7").Value = "GooD d AsomE Mobile, Battery  was so gooD and fast charge   i love this phonE camera was so nicE"
$ws.Range("B8").Value = "Good value for money and a phone with good camera and great battery life with fast charging which is a life saver sometimes"
$ws.Range("B9").Value = "Best mobile ...very good camera quality ... Fast refresh rate... Value of money.. must enjoyed this mobile"
$ws.Range("B10").Value = "I overall like this phone. Good performance, camera, battery is also better not so high. Still I am satisfied with the phone"
$ws.Range("B11").Value = "Very awesome mobile. Camera is good but front camera is not that much better. Overall it is better than samsung M21 or M31 if you are planning for these smartphones so change your desigen and buy Realme 7. Charging is tooo good."
